$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has an unused/duplicate "_1" column (C) sitting between the
# "South Africa price price" (B) and "South Africa points points" (D)
# statistics columns. Remove it so the points column shifts left into C,
# and clean up the duplicated-word header text for both stat columns.
$ws.Columns.Item(3).Delete()

$ws.Range("B1").Value = "South Africa_priceprice"
$ws.Range("C1").Value = "South Africa_pointspoints"

# Restore the intended column widths (price column slightly narrower,
# points column slightly wider) now that the spacer column is gone.
$ws.Columns.Item(2).ColumnWidth = 24.8776041666667
$ws.Columns.Item(3).ColumnWidth = 26.8776041666667
